$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row before the current row 27, pushing the existing
# rows 27-42 down to 28-43 (dimension grows from A1:T42 to A1:T43).
$ws.Rows(27).Insert()

# Populate the newly inserted row 27 with the new record's data.
$ws.Cells.Item(27, 1).Value = 6
$ws.Cells.Item(27, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(27, 3).Value = "Metropolitana"
$ws.Cells.Item(27, 4).Value = 44813
$ws.Cells.Item(27, 5).Value = 13
$ws.Cells.Item(27, 6).Value = "Fruta"
$ws.Cells.Item(27, 7).Value = 100102
$ws.Cells.Item(27, 8).Value = "Cítricos"
$ws.Cells.Item(27, 9).Value = 100102006
$ws.Cells.Item(27, 10).Value = "Pomelo"
$ws.Cells.Item(27, 11).Value = "Start Ruby"
$ws.Cells.Item(27, 12).Value = "Primera"
$ws.Cells.Item(27, 13).Value = 20
$ws.Cells.Item(27, 14).Value = 160000
$ws.Cells.Item(27, 15).Value = 160000
$ws.Cells.Item(27, 16).Value = 160000
$ws.Cells.Item(27, 17).Value = "`$/bins (350 kilos)"
$ws.Cells.Item(27, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(27, 19).Value = 457
$ws.Cells.Item(27, 20).Value = 350
